$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The rows 136-142 (columns A:F) need to be reordered (a pure row permutation,
# reverting an earlier reorder of the underlying mapp_differences export).
# new row <- old row
#   136 <- 139
#   137 <- 142
#   138 <- 141
#   139 <- 138
#   140 <- 137
#   141 <- 136
#   142 <- 140
#
# Use Range.Copy (not .Value = "...") so the original shared-string /
# text typing of each cell (IDs, bracket messages, epoch timestamps stored
# as text, and date strings) is preserved exactly instead of being
# re-parsed (which would turn numeric-looking text like "1727957259" into
# a real number).

# Stage the current rows 136-142 far away first, since several of the
# source/destination rows overlap (this is a permutation, not a simple
# shift), so a direct cell-to-cell copy could overwrite source data
# before it has been read.
$ws.Range("A136:F142").Copy($ws.Range("A1000:F1006"))

# Staging row offsets (1000 = old row 136, 1001 = old row 137, ... 1006 = old row 142)
function OldRow([int]$oldRow) { return 1000 + ($oldRow - 136) }

$ws.Range("A" + (OldRow 139) + ":F" + (OldRow 139)).Copy($ws.Range("A136:F136"))
$ws.Range("A" + (OldRow 142) + ":F" + (OldRow 142)).Copy($ws.Range("A137:F137"))
$ws.Range("A" + (OldRow 141) + ":F" + (OldRow 141)).Copy($ws.Range("A138:F138"))
$ws.Range("A" + (OldRow 138) + ":F" + (OldRow 138)).Copy($ws.Range("A139:F139"))
$ws.Range("A" + (OldRow 137) + ":F" + (OldRow 137)).Copy($ws.Range("A140:F140"))
$ws.Range("A" + (OldRow 136) + ":F" + (OldRow 136)).Copy($ws.Range("A141:F141"))
$ws.Range("A" + (OldRow 140) + ":F" + (OldRow 140)).Copy($ws.Range("A142:F142"))

# Clean up the staging area.
$ws.Range("A1000:F1006").ClearContents()
